$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing counters
$ws.Range("D2").Value = 30
$ws.Range("D9").Value = 21
$ws.Range("D10").Value = 30

# Duplicate row 10 (values + formatting) into row 11 so the new row inherits
# the same text/shared-string typing and styles as the existing "asd" row,
# then overwrite the cells that differ for the new "poooq" user.
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F11").PasteSpecial(-4104)  # xlPasteAll

$ws.Range("A11").Value = "poooq"
$ws.Range("B11").Value = "poooq"
$ws.Range("C11").Value = "E9956AF6"
$ws.Range("D11").Value = 2
